$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.157.96'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.911.68'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.81%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7389'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '245.81'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3105'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.65'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06995'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08049'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7708'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.895.44'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.364'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.31'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.45'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.166.84'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.987'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007878'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.61'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.92%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.06%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.147.13'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.174'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +7.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.451'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.71'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.02'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1285'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.055'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.41%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.560'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.49%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.354'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.95%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.350'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.088'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.85%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.317'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.18%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05171'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7526'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.733'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01958'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.795'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.362'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4532'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.87'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.990'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.966'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +5.48%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.002'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.23%  '

$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.8418'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.02'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.87'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.60%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.05'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.12%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.043.51'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.64%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1189'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.38%  '
